$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet (JPMGE), matching the
# position/sheetId/rId sequence Excel produces for an appended sheet.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "TwoxTwoCET-Scalar"

$ws.Range("A2").Value = "'A.L"
$ws.Range("A3").Value = "'B.L"
$ws.Range("A4").Value = "'W.L"
$ws.Range("A5").Value = "'PX.L"
$ws.Range("A6").Value = "'PY.L"
$ws.Range("A7").Value = "'PW.L"
$ws.Range("A8").Value = "'PL.L"
$ws.Range("A9").Value = "'PK.L"
$ws.Range("A10").Value = "'CONS.L"
$ws.Range("A11").Value = "'SAX.L"
$ws.Range("A12").Value = "'SAY.L"
$ws.Range("A13").Value = "'SBX.L"
$ws.Range("A14").Value = "'SBY.L"
$ws.Range("A15").Value = "'DAL.L"
$ws.Range("A16").Value = "'DAK.L"
$ws.Range("A17").Value = "'DBL.L"
$ws.Range("A18").Value = "'DBK.L"
$ws.Range("A19").Value = "'SW.L"
$ws.Range("A20").Value = "'DWX.L"
$ws.Range("A21").Value = "'DWY.L"
$ws.Range("A22").Value = "'DW.L"
$ws.Range("A23").Value = "'CWI.L"
$ws.Range("A24").Value = "'PX.L/PW.L"
$ws.Range("A25").Value = "'PY.L/PW.L"
$ws.Range("A26").Value = "'PW.L/PW.L"
$ws.Range("A27").Value = "'PL.L/PW.L"
$ws.Range("A28").Value = "'PK.L/PW.L"
$ws.Range("A29").Value = "'CONS.L/PW.L"
$ws.Range("B1").Value = "'benchmark"
$ws.Range("C1").Value = "'diff=10"
$ws.Range("D1").Value = "'PW.FX=1"
$ws.Range("E1").Value = "'TrA=2, TrB=1.5"
$ws.Range("F1").Value = "'TA=0.1"
$ws.Range("G1").Value = "'TA=100%"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.84841617526015112
$ws.Range("D2").Value = 0.8484161751970315
$ws.Range("E2").Value = 0.70716515735515006
$ws.Range("F2").Value = 0.38751484805497743
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1.1506325096039884
$ws.Range("D3").Value = 1.1506325095207364
$ws.Range("E3").Value = 1.2893059702343606
$ws.Range("F3").Value = 1.5972629436210768
$ws.Range("G3").Value = 1.9601317042077897
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.0035850931326258
$ws.Range("D4").Value = 1.0035850930478027
$ws.Range("E4").Value = 1.0068914996667162
$ws.Range("F4").Value = 0.99916851890038427
$ws.Range("G4").Value = 0.97073779227082346
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1.0231427182799986
$ws.Range("D5").Value = 1.026810780303149
$ws.Range("E5").Value = 1.0161855054344244
$ws.Range("F5").Value = 1.0862618717772732
$ws.Range("G5").Value = 1.1690468695684746
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0.87359232216492444
$ws.Range("D6").Value = 0.87672423219236795
$ws.Range("E6").Value = 0.88681325620491058
$ws.Range("F6").Value = 0.8240896931384919
$ws.Range("G6").Value = 0.76012909543532714
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0.99642771370717265
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1.0314650721357199
$ws.Range("D8").Value = 1.0351629708854646
$ws.Range("E8").Value = 1.0678356276899741
$ws.Range("F8").Value = 1.1033933594907179
$ws.Range("G8").Value = 1.1648853507249877
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0.96853492724769141
$ws.Range("D9").Value = 0.97200721521014155
$ws.Range("E9").Value = 0.94594737164343623
$ws.Range("F9").Value = 0.85817103861968913
$ws.Range("G9").Value = 0.77659023381665881
$ws.Range("B10").Value = 200
$ws.Range("C10").Value = 200
$ws.Range("D10").Value = 200.71701860956063
$ws.Range("E10").Value = 201.37829993334105
$ws.Range("F10").Value = 199.83370377972741
$ws.Range("G10").Value = 194.14755845416465
$ws.Range("B11").Value = 80
$ws.Range("C11").Value = 80.000000000000014
$ws.Range("D11").Value = 80
$ws.Range("E11").Value = 83.789979418307453
$ws.Range("F11").Value = 86.636474448757937
$ws.Range("G11").Value = "'Undf"
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = 20.000000000000004
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 15.953297102016906
$ws.Range("F12").Value = 12.465814866890884
$ws.Range("G12").Value = "'Undf"
$ws.Range("B13").Value = 20
$ws.Range("C13").Value = 30.000000000000004
$ws.Range("D13").Value = 30
$ws.Range("E13").Value = 34.553784903611799
$ws.Range("F13").Value = 39.310669882974452
$ws.Range("G13").Value = 44.380086793200221
$ws.Range("B14").Value = 80
$ws.Range("C14").Value = 80
$ws.Range("D14").Value = 80
$ws.Range("E14").Value = 75.119576065481453
$ws.Range("F14").Value = 69.26902850423329
$ws.Range("G14").Value = 62.049748999369605
$ws.Range("B15").Value = 40
$ws.Range("C15").Value = 38.517354230045243
$ws.Range("D15").Value = 38.517354222141414
$ws.Range("E15").Value = 37.194397849630214
$ws.Range("F15").Value = 34.40060528627518
$ws.Range("G15").Value = "'Undf"
$ws.Range("B16").Value = 60
$ws.Range("C16").Value = 61.530004404080294
$ws.Range("D16").Value = 61.53000441249764
$ws.Range("E16").Value = 62.98051725431754
$ws.Range("F16").Value = 66.345864158482783
$ws.Range("G16").Value = "'Undf"
$ws.Range("B17").Value = 60
$ws.Range("C17").Value = 58.508040668387657
$ws.Range("D17").Value = 58.508040660383685
$ws.Range("E17").Value = 57.160534034089842
$ws.Range("F17").Value = 54.261106485862463
$ws.Range("G17").Value = 51.016980025031643
$ws.Range("B18").Value = 40
$ws.Range("C18").Value = 41.539717147859804
$ws.Range("D18").Value = 41.539717156383809
$ws.Range("E18").Value = 43.017230886987115
$ws.Range("F18").Value = 46.51080952457405
$ws.Range("G18").Value = 51.016980025031621
$ws.Range("B19").Value = 200
$ws.Range("C19").Value = 200.71701862652515
$ws.Range("D19").Value = 200.71701860956054
$ws.Range("E19").Value = 201.37829993334324
$ws.Range("F19").Value = 199.83370378007686
$ws.Range("G19").Value = 194.14755845416468
$ws.Range("B20").Value = 100
$ws.Range("C20").Value = 102.3922693035888
$ws.Range("D20").Value = 102.39226930136157
$ws.Range("E20").Value = 103.80375515055297
$ws.Range("F20").Value = 96.362396524949176
$ws.Range("G20").Value = 86.990815158845137
$ws.Range("B21").Value = 100
$ws.Range("C21").Value = 109.01892428236799
$ws.Range("D21").Value = 109.0189242656265
$ws.Range("E21").Value = 108.13373375818081
$ws.Range("F21").Value = 115.47154072457153
$ws.Range("G21").Value = 121.62568025179993
$ws.Range("B22").Value = 200
$ws.Range("C22").Value = 200.71701859259665
$ws.Range("D22").Value = 200.71701860956063
$ws.Range("E22").Value = 201.37829993334105
$ws.Range("F22").Value = 199.83370377972741
$ws.Range("G22").Value = 194.14755845416465
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1.0035850929629833
$ws.Range("D23").Value = 1.0035850930478032
$ws.Range("E23").Value = 1.0068914996667053
$ws.Range("F23").Value = 0.99916851889863711
$ws.Range("G23").Value = 0.97073779227082324
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 1.0268107803559916
$ws.Range("D24").Value = 1.026810780303149
$ws.Range("E24").Value = 1.0161855054344244
$ws.Range("F24").Value = 1.0862618717772732
$ws.Range("G24").Value = 1.1690468695684746
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 0.87672423212192319
$ws.Range("D25").Value = 0.87672423219236795
$ws.Range("E25").Value = 0.88681325620491058
$ws.Range("F25").Value = 0.8240896931384919
$ws.Range("G25").Value = 0.76012909543532714
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 1.0351629706265315
$ws.Range("D27").Value = 1.0351629708854646
$ws.Range("E27").Value = 1.0678356276899741
$ws.Range("F27").Value = 1.1033933594907179
$ws.Range("G27").Value = 1.1648853507249877
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 0.9720072152994349
$ws.Range("D28").Value = 0.97200721521014155
$ws.Range("E28").Value = 0.94594737164343623
$ws.Range("F28").Value = 0.85817103861968913
$ws.Range("G28").Value = 0.77659023381665881
$ws.Range("B29").Value = 200
$ws.Range("C29").Value = 200.71701865447656
$ws.Range("D29").Value = 200.71701860956063
$ws.Range("E29").Value = 201.37829993334105
$ws.Range("F29").Value = 199.83370377972741
$ws.Range("G29").Value = 194.14755845416465

